$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 12 (USER_GEN_REC -> USER_REC) ---
$ws.Range("A12").Value = "USER_REC"
$ws.Range("B12").Value = "Test recipe recogniser  when you upload an image"
$ws.Range("C12").Value = "Successful recipe recogniser "

# --- Clear rows 13 & 14 (former USER_FIL_PREF / USER_INSUFF_ERR rows) ---
$ws.Range("A13").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# --- Give the whole column B the same body font as the rest of the sheet ---
# (keeps each row's own wrap setting untouched)
$ws.Columns.Item(2).Font.Name = "Aptos Narrow"
$ws.Columns.Item(2).Font.Size = 11

# Rows that should keep their wrapped, taller presentation
$ws.Range("B13").WrapText = $true
$ws.Range("B14").WrapText = $true

# Row heights: row 11 loses its explicit 16pt height, row 12 gains it
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 16
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15

# --- Selection moves to B13 ---
$ws.Range("B13").Select()
